$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.332.32"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "2.246.88"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'230.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").Value = "'64.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("D10").Value = "'0.0948"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.79%  "
$ws.Range("D11").Value = "'56.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "'26.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "2.581.79"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "'14.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.34%  "
$ws.Range("D16").Value = "'6.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "2.255.32"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "43.210.39"
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("E20").Value = "  -5.20%  "
$ws.Range("D21").Value = "'72.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").Value = "'3.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.90%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'173.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.92%  "
$ws.Range("D30").Value = "'21.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.50%  "
$ws.Range("D31").Value = "'1.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.82%  "
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").Value = "'0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").Value = "'3.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("D38").Value = "'6.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("E40").Value = "  -3.30%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'8.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.32%  "
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").Value = "'17.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.06%  "
$ws.Range("D45").Value = "'10.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "'96.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "'1.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0935"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("B49").Value = "TerraClassic"
$ws.Range("C49").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D49").Value = "'0.000205"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.423.60"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("E51").Value = "  -1.90%  "
